$wb = $excel.ActiveWorkbook

# Sheet "展览": increment the "想去人数" (F column) for several rows by 1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7147
$ws1.Range("F6").Value = 559
$ws1.Range("F7").Value = 174
$ws1.Range("F18").Value = 4
$ws1.Range("F27").Value = 285
$ws1.Range("F37").Value = 129

# Sheet "全部类型": same rows (offset by one due to extra 演出 row) get the same increments
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7147
$ws4.Range("F7").Value = 559
$ws4.Range("F8").Value = 174
$ws4.Range("F19").Value = 4
$ws4.Range("F28").Value = 285
$ws4.Range("F38").Value = 129
